$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 9093823
$ws.Cells.Item(17, 10).Value = 9093823
$ws.Cells.Item(17, 12).Value = 27281469
$ws.Cells.Item(17, 14).Value = -27281805
$ws.Cells.Item(28, 8).Value = 411.93332
$ws.Cells.Item(28, 9).Value = 119.92857
$ws.Cells.Item(28, 11).Value = 119.92857
$ws.Cells.Item(28, 13).Value = 365.07143
$ws.Cells.Item(33, 8).Value = 669.3333
$ws.Cells.Item(33, 9).Value = 574.1177
$ws.Cells.Item(33, 11).Value = 574.1177
$ws.Cells.Item(33, 13).Value = -345.1177
$ws.Cells.Item(54, 8).Value = 2073
$ws.Cells.Item(54, 9).Value = 2076
$ws.Cells.Item(54, 10).Value = 2070
$ws.Cells.Item(54, 11).Value = 2076
$ws.Cells.Item(54, 12).Value = 2070
$ws.Cells.Item(54, 13).Value = -1590
$ws.Cells.Item(54, 14).Value = -3042
$ws.Cells.Item(98, 8).Value = 609.913
$ws.Cells.Item(98, 9).Value = 525.1429000000001
$ws.Cells.Item(98, 11).Value = 525.1429000000001
$ws.Cells.Item(98, 13).Value = 972.8570999999999
$ws.Cells.Item(107, 8).Value = 412.44446
$ws.Cells.Item(107, 9).Value = 424.35715
$ws.Cells.Item(107, 10).Value = 370.75
$ws.Cells.Item(107, 11).Value = 424.35715
$ws.Cells.Item(107, 12).Value = 370.75
$ws.Cells.Item(107, 13).Value = 1495.64285
$ws.Cells.Item(107, 14).Value = -4210.75
$ws.Cells.Item(112, 8).Value = 1497.3334
$ws.Cells.Item(112, 9).Value = 1296.25
$ws.Cells.Item(112, 10).Value = 1899.5
$ws.Cells.Item(112, 11).Value = 3888.75
$ws.Cells.Item(112, 12).Value = 5698.5
$ws.Cells.Item(112, 13).Value = -2780.75
$ws.Cells.Item(112, 14).Value = -7914.5
$ws.Cells.Item(122, 8).Value = 609.913
$ws.Cells.Item(122, 9).Value = 525.1429000000001
$ws.Cells.Item(122, 11).Value = 1575.4287
$ws.Cells.Item(122, 13).Value = 874.5712999999998
$ws.Cells.Item(138, 8).Value = 3563.5227
$ws.Cells.Item(138, 10).Value = 3178.5483
$ws.Cells.Item(138, 12).Value = 9535.644899999999
$ws.Cells.Item(138, 14).Value = -19815.6449

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1163.8064
$ws.Cells.Item(2, 9).Value = 946.5185
$ws.Cells.Item(2, 11).Value = 946.5185
$ws.Cells.Item(2, 13).Value = -833.5185
$ws.Cells.Item(10, 8).Value = 8000
$ws.Cells.Item(10, 10).Value = 8000
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 14).Value = -8340
$ws.Cells.Item(74, 8).Value = 4531.433
$ws.Cells.Item(74, 9).Value = 3694.1538
$ws.Cells.Item(74, 10).Value = 9973.75
$ws.Cells.Item(74, 11).Value = 3694.1538
$ws.Cells.Item(74, 12).Value = 9973.75
$ws.Cells.Item(74, 13).Value = -2820.1538
$ws.Cells.Item(74, 14).Value = -11721.75
$ws.Cells.Item(77, 8).Value = 4531.433
$ws.Cells.Item(77, 9).Value = 3694.1538
$ws.Cells.Item(77, 10).Value = 9973.75
$ws.Cells.Item(77, 11).Value = 18470.769
$ws.Cells.Item(77, 12).Value = 49868.75
$ws.Cells.Item(77, 13).Value = -14102.769
$ws.Cells.Item(77, 14).Value = -58604.75
$ws.Cells.Item(102, 8).Value = 4767.6313
$ws.Cells.Item(102, 9).Value = 2839.5334
$ws.Cells.Item(102, 11).Value = 2839.5334
$ws.Cells.Item(102, 13).Value = -1217.5334
$ws.Cells.Item(116, 8).Value = 1163.8064
$ws.Cells.Item(116, 9).Value = 946.5185
$ws.Cells.Item(116, 11).Value = 946.5185
$ws.Cells.Item(116, 13).Value = 1347.4815
$ws.Cells.Item(122, 8).Value = 2104.4
$ws.Cells.Item(122, 9).Value = 2209.8
$ws.Cells.Item(122, 11).Value = 6629.400000000001
$ws.Cells.Item(122, 13).Value = -4179.400000000001
$ws.Cells.Item(139, 8).Value = 92666.5
$ws.Cells.Item(139, 10).Value = 92666.5
$ws.Cells.Item(139, 12).Value = 92666.5
$ws.Cells.Item(139, 14).Value = -102946.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1163.8064
$ws.Cells.Item(3, 9).Value = 946.5185
$ws.Cells.Item(3, 11).Value = 946.5185
$ws.Cells.Item(3, 13).Value = -832.5185
$ws.Cells.Item(20, 8).Value = 3927.5518
$ws.Cells.Item(20, 9).Value = 2997.1
$ws.Cells.Item(20, 11).Value = 2997.1
$ws.Cells.Item(20, 13).Value = -2750.1
$ws.Cells.Item(94, 8).Value = 1290.6
$ws.Cells.Item(94, 9).Value = 877.3333
$ws.Cells.Item(94, 10).Value = 5010
$ws.Cells.Item(94, 11).Value = 877.3333
$ws.Cells.Item(94, 12).Value = 5010
$ws.Cells.Item(94, 13).Value = -426.3333
$ws.Cells.Item(94, 14).Value = -5912
$ws.Cells.Item(137, 8).Value = 74999
$ws.Cells.Item(137, 10).Value = 74999
$ws.Cells.Item(137, 12).Value = 74999
$ws.Cells.Item(137, 14).Value = -85199

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 100000
$ws.Cells.Item(52, 10).Value = 100000
$ws.Cells.Item(52, 12).Value = 100000
$ws.Cells.Item(52, 14).Value = -100588
$ws.Cells.Item(97, 8).Value = 25266
$ws.Cells.Item(97, 10).Value = 26022.572
$ws.Cells.Item(97, 12).Value = 26022.572
$ws.Cells.Item(97, 14).Value = -28004.572
$ws.Cells.Item(122, 8).Value = 3040
$ws.Cells.Item(122, 9).Value = 3344.1765
$ws.Cells.Item(122, 11).Value = 10032.5295
$ws.Cells.Item(122, 13).Value = -7582.529500000001
$ws.Cells.Item(135, 8).Value = 99429.625
$ws.Cells.Item(135, 10).Value = 100776.71
$ws.Cells.Item(135, 12).Value = 100776.71
$ws.Cells.Item(135, 14).Value = -110916.71
$ws.Cells.Item(137, 8).Value = 40700
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 198250
$ws.Cells.Item(139, 10).Value = 198250
$ws.Cells.Item(139, 12).Value = 198250
$ws.Cells.Item(139, 14).Value = -208530

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 354.8
$ws.Cells.Item(33, 9).Value = 305.4
$ws.Cells.Item(33, 10).Value = 379.5
$ws.Cells.Item(33, 11).Value = 1832.4
$ws.Cells.Item(33, 12).Value = 2277
$ws.Cells.Item(33, 13).Value = -1549.4
$ws.Cells.Item(33, 14).Value = -2843
$ws.Cells.Item(120, 8).Value = 13088.333
$ws.Cells.Item(120, 9).Value = 6176.6665
$ws.Cells.Item(120, 11).Value = 18529.9995
$ws.Cells.Item(120, 13).Value = -13691.9995
$ws.Cells.Item(139, 8).Value = 26938.8
$ws.Cells.Item(139, 9).Value = 8673.75
$ws.Cells.Item(139, 11).Value = 26021.25
$ws.Cells.Item(139, 13).Value = -20881.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 28190
$ws.Cells.Item(26, 10).Value = 28987.5
$ws.Cells.Item(26, 12).Value = 28987.5
$ws.Cells.Item(26, 14).Value = -29547.5
$ws.Cells.Item(50, 8).Value = 28190
$ws.Cells.Item(50, 10).Value = 28987.5
$ws.Cells.Item(50, 12).Value = 28987.5
$ws.Cells.Item(50, 14).Value = -29983.5
$ws.Cells.Item(80, 8).Value = 2624.5715
$ws.Cells.Item(80, 9).Value = 2500
$ws.Cells.Item(80, 10).Value = 2718
$ws.Cells.Item(80, 11).Value = 2500
$ws.Cells.Item(80, 12).Value = 2718
$ws.Cells.Item(80, 13).Value = -1502
$ws.Cells.Item(80, 14).Value = -4714
$ws.Cells.Item(83, 8).Value = 2624.5715
$ws.Cells.Item(83, 9).Value = 2500
$ws.Cells.Item(83, 10).Value = 2718
$ws.Cells.Item(83, 11).Value = 12500
$ws.Cells.Item(83, 12).Value = 13590
$ws.Cells.Item(83, 13).Value = -7508
$ws.Cells.Item(83, 14).Value = -23574
$ws.Cells.Item(97, 8).Value = 741.5862
$ws.Cells.Item(97, 9).Value = 759.3182
$ws.Cells.Item(97, 11).Value = 759.3182
$ws.Cells.Item(97, 13).Value = -263.3182
$ws.Cells.Item(122, 8).Value = 1633.96
$ws.Cells.Item(122, 9).Value = 1583.2858
$ws.Cells.Item(122, 10).Value = 1900
$ws.Cells.Item(122, 11).Value = 4749.857400000001
$ws.Cells.Item(122, 12).Value = 5700
$ws.Cells.Item(122, 13).Value = -2299.857400000001
$ws.Cells.Item(122, 14).Value = -10600

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3179.5
$ws.Cells.Item(7, 9).Value = 3179.5
$ws.Cells.Item(7, 11).Value = 3179.5
$ws.Cells.Item(7, 13).Value = -3067.5
$ws.Cells.Item(22, 8).Value = 2399
$ws.Cells.Item(22, 9).Value = 2399
$ws.Cells.Item(22, 11).Value = 2399
$ws.Cells.Item(22, 13).Value = -2104
$ws.Cells.Item(27, 8).Value = 2399
$ws.Cells.Item(27, 9).Value = 2399
$ws.Cells.Item(27, 11).Value = 2399
$ws.Cells.Item(27, 13).Value = -2292
$ws.Cells.Item(46, 8).Value = 3159.2
$ws.Cells.Item(46, 9).Value = 3449
$ws.Cells.Item(46, 10).Value = 2000
$ws.Cells.Item(46, 11).Value = 3449
$ws.Cells.Item(46, 12).Value = 2000
$ws.Cells.Item(46, 13).Value = -3261
$ws.Cells.Item(46, 14).Value = -2376
$ws.Cells.Item(82, 8).Value = 1466
$ws.Cells.Item(82, 9).Value = 1474.5
$ws.Cells.Item(82, 10).Value = 1449
$ws.Cells.Item(82, 11).Value = 1474.5
$ws.Cells.Item(82, 12).Value = 1449
$ws.Cells.Item(82, 13).Value = -1113.5
$ws.Cells.Item(82, 14).Value = -2171
$ws.Cells.Item(85, 8).Value = 1466
$ws.Cells.Item(85, 9).Value = 1474.5
$ws.Cells.Item(85, 10).Value = 1449
$ws.Cells.Item(85, 11).Value = 1474.5
$ws.Cells.Item(85, 12).Value = 1449
$ws.Cells.Item(85, 13).Value = -226.5
$ws.Cells.Item(85, 14).Value = -3945
$ws.Cells.Item(122, 8).Value = 5006
$ws.Cells.Item(122, 9).Value = 3852.4285
$ws.Cells.Item(122, 10).Value = 10389.333
$ws.Cells.Item(122, 11).Value = 11557.2855
$ws.Cells.Item(122, 12).Value = 31167.999
$ws.Cells.Item(122, 13).Value = -9107.2855
$ws.Cells.Item(122, 14).Value = -36067.999
$ws.Cells.Item(126, 8).Value = 3179.5
$ws.Cells.Item(126, 9).Value = 3179.5
$ws.Cells.Item(126, 11).Value = 9538.5
$ws.Cells.Item(126, 13).Value = -7068.5
$ws.Cells.Item(132, 8).Value = 8454.52
$ws.Cells.Item(132, 9).Value = 8832.257
$ws.Cells.Item(132, 11).Value = 26496.771
$ws.Cells.Item(132, 13).Value = -23966.771

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 5000
$ws.Cells.Item(38, 10).Value = 5000
$ws.Cells.Item(38, 12).Value = 5000
$ws.Cells.Item(38, 14).Value = -5946
$ws.Cells.Item(99, 8).Value = 40580
$ws.Cells.Item(99, 10).Value = 40580
$ws.Cells.Item(99, 12).Value = 40580
$ws.Cells.Item(99, 14).Value = -46570
$ws.Cells.Item(122, 8).Value = 3287.9814
$ws.Cells.Item(122, 9).Value = 2532.1462
$ws.Cells.Item(122, 11).Value = 7596.4386
$ws.Cells.Item(122, 13).Value = -5146.4386
$ws.Cells.Item(132, 8).Value = 3171.459
$ws.Cells.Item(132, 9).Value = 2837.745
$ws.Cells.Item(132, 10).Value = 4873.4
$ws.Cells.Item(132, 11).Value = 8513.235000000001
$ws.Cells.Item(132, 12).Value = 14620.2
$ws.Cells.Item(132, 13).Value = -5983.235000000001
$ws.Cells.Item(132, 14).Value = -19680.2

